$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Dezembro")
foreach ($r in @(3,6,9,10,11,13,14,17,18)) {
  $ws2.Rows.Item($r).OutlineLevel = 0
}
Write-Host "done"
